$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value2 = 2138.9722
$ws.Range("I15").Value2 = 2138.9722
$ws.Range("K15").Value2 = 6416.9166
$ws.Range("M15").Value2 = -6247.9166
$ws.Range("H125").Value2 = 1787.4546
$ws.Range("I125").Value2 = 2666.6667
$ws.Range("J125").Value2 = 1457.75
$ws.Range("K125").Value2 = 24000.0003
$ws.Range("L125").Value2 = 13119.75
$ws.Range("M125").Value2 = -21540.0003
$ws.Range("N125").Value2 = -18039.75

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value2 = 100
$ws.Range("I4").Value2 = 100
$ws.Range("J4").Value2 = 0
$ws.Range("K4").Value2 = 100
$ws.Range("L4").Value2 = 0
$ws.Range("M4").Value2 = 16
$ws.Range("N4").ClearContents()
$ws.Range("H32").Value2 = 19777.97
$ws.Range("I32").Value2 = 17230.393
$ws.Range("K32").Value2 = 17230.393
$ws.Range("M32").Value2 = -16943.393
$ws.Range("H61").Value2 = 2966.9736
$ws.Range("I61").Value2 = 2336.739
$ws.Range("J61").Value2 = 3933.3333
$ws.Range("K61").Value2 = 2336.739
$ws.Range("L61").Value2 = 3933.3333
$ws.Range("M61").Value2 = -2124.739
$ws.Range("N61").Value2 = -4357.3333
$ws.Range("H74").Value2 = 1589.0571
$ws.Range("I74").Value2 = 1153.9333
$ws.Range("J74").Value2 = 4199.8
$ws.Range("K74").Value2 = 1153.9333
$ws.Range("L74").Value2 = 4199.8
$ws.Range("M74").Value2 = -279.9332999999999
$ws.Range("N74").Value2 = -5947.8
$ws.Range("H77").Value2 = 1589.0571
$ws.Range("I77").Value2 = 1153.9333
$ws.Range("J77").Value2 = 4199.8
$ws.Range("K77").Value2 = 5769.666499999999
$ws.Range("L77").Value2 = 20999
$ws.Range("M77").Value2 = -1401.666499999999
$ws.Range("N77").Value2 = -29735
$ws.Range("H109").Value2 = 26666.666
$ws.Range("J109").Value2 = 26666.666
$ws.Range("L109").Value2 = 26666.666
$ws.Range("N109").Value2 = -29440.666
$ws.Range("H122").Value2 = 3586.6428
$ws.Range("I122").Value2 = 2479.8
$ws.Range("J122").Value2 = 4201.5557
$ws.Range("K122").Value2 = 7439.400000000001
$ws.Range("L122").Value2 = 12604.6671
$ws.Range("M122").Value2 = -4989.400000000001
$ws.Range("N122").Value2 = -17504.6671
$ws.Range("H124").Value2 = 0
$ws.Range("J124").Value2 = 0
$ws.Range("L124").Value2 = 0
$ws.Range("N124").ClearContents()
$ws.Range("H128").Value2 = 0
$ws.Range("J128").Value2 = 0
$ws.Range("L128").Value2 = 0
$ws.Range("N128").ClearContents()
$ws.Range("H129").Value2 = 0
$ws.Range("J129").Value2 = 0
$ws.Range("L129").Value2 = 0
$ws.Range("N129").ClearContents()
$ws.Range("H131").Value2 = 0
$ws.Range("J131").Value2 = 0
$ws.Range("L131").Value2 = 0
$ws.Range("N131").ClearContents()
$ws.Range("H132").Value2 = 1884.4728
$ws.Range("I132").Value2 = 1285.9269
$ws.Range("K132").Value2 = 3857.7807
$ws.Range("M132").Value2 = -1327.7807
$ws.Range("H136").Value2 = 2966.9736
$ws.Range("I136").Value2 = 2336.739
$ws.Range("J136").Value2 = 3933.3333
$ws.Range("K136").Value2 = 7010.217000000001
$ws.Range("L136").Value2 = 11799.9999
$ws.Range("M136").Value2 = -4460.217000000001
$ws.Range("N136").Value2 = -16899.9999

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H53").Value2 = 0
$ws.Range("J53").Value2 = 0
$ws.Range("L53").Value2 = 0
$ws.Range("N53").ClearContents()
$ws.Range("H86").Value2 = 27374
$ws.Range("I86").Value2 = 2480
$ws.Range("J86").Value2 = 52268
$ws.Range("K86").Value2 = 2480
$ws.Range("L86").Value2 = 52268
$ws.Range("M86").Value2 = -1357
$ws.Range("N86").Value2 = -54514
$ws.Range("H89").Value2 = 27374
$ws.Range("I89").Value2 = 2480
$ws.Range("J89").Value2 = 52268
$ws.Range("K89").Value2 = 12400
$ws.Range("L89").Value2 = 261340
$ws.Range("M89").Value2 = -6784
$ws.Range("N89").Value2 = -272572
$ws.Range("H102").Value2 = 0
$ws.Range("I102").Value2 = 0
$ws.Range("K102").Value2 = 0
$ws.Range("M102").ClearContents()
$ws.Range("H134").Value2 = 3166.2104
$ws.Range("I134").Value2 = 3026.2646
$ws.Range("J134").Value2 = 4355.75
$ws.Range("K134").Value2 = 9078.7938
$ws.Range("L134").Value2 = 13067.25
$ws.Range("M134").Value2 = -6543.793799999999
$ws.Range("N134").Value2 = -18137.25

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value2 = 2226285
$ws.Range("I31").Value2 = 3573676.5
$ws.Range("J31").Value2 = 7052.2354
$ws.Range("K31").Value2 = 3573676.5
$ws.Range("L31").Value2 = 7052.2354
$ws.Range("M31").Value2 = -3573381.5
$ws.Range("N31").Value2 = -7642.2354
$ws.Range("H34").Value2 = 2226285
$ws.Range("I34").Value2 = 3573676.5
$ws.Range("J34").Value2 = 7052.2354
$ws.Range("K34").Value2 = 3573676.5
$ws.Range("L34").Value2 = 7052.2354
$ws.Range("M34").Value2 = -3573474.5
$ws.Range("N34").Value2 = -7456.2354
$ws.Range("H99").Value2 = 6487.6665
$ws.Range("I99").Value2 = 4978
$ws.Range("K99").Value2 = 4978
$ws.Range("M99").Value2 = -3480
$ws.Range("H126").Value2 = 6487.6665
$ws.Range("I126").Value2 = 4978
$ws.Range("K126").Value2 = 14934
$ws.Range("M126").Value2 = -12464

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value2 = 1307.6154
$ws.Range("I114").Value2 = 633.3333
$ws.Range("J114").Value2 = 1885.5714
$ws.Range("K114").Value2 = 1899.9999
$ws.Range("L114").Value2 = 5656.7142
$ws.Range("M114").Value2 = 1354.0001
$ws.Range("N114").Value2 = -12164.7142
$ws.Range("H132").Value2 = 2499.9312
$ws.Range("I132").Value2 = 1373.4546
$ws.Range("J132").Value2 = 3188.3333
$ws.Range("K132").Value2 = 12361.0914
$ws.Range("L132").Value2 = 28694.9997
$ws.Range("M132").Value2 = -9831.091400000001
$ws.Range("N132").Value2 = -33754.9997

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value2 = 4294.6665
$ws.Range("I132").Value2 = 4428.769
$ws.Range("J132").Value2 = 3946
$ws.Range("K132").Value2 = 13286.307
$ws.Range("L132").Value2 = 11838
$ws.Range("M132").Value2 = -10756.307
$ws.Range("N132").Value2 = -16898

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value2 = 17666.666
$ws.Range("I40").Value2 = 21500
$ws.Range("J40").Value2 = 10000
$ws.Range("K40").Value2 = 21500
$ws.Range("L40").Value2 = 10000
$ws.Range("M40").Value2 = -21364
$ws.Range("N40").Value2 = -10272
$ws.Range("H100").Value2 = 3394
$ws.Range("I100").Value2 = 1730
$ws.Range("J100").Value2 = 3810
$ws.Range("K100").Value2 = 1730
$ws.Range("L100").Value2 = 3810
$ws.Range("M100").Value2 = -1189
$ws.Range("N100").Value2 = -4892
$ws.Range("H122").Value2 = 6583
$ws.Range("I122").Value2 = 5999.5
$ws.Range("J122").Value2 = 6874.75
$ws.Range("K122").Value2 = 17998.5
$ws.Range("L122").Value2 = 20624.25
$ws.Range("M122").Value2 = -15548.5
$ws.Range("N122").Value2 = -25524.25
$ws.Range("H132").Value2 = 3319.6365
$ws.Range("I132").Value2 = 2534.25
$ws.Range("J132").Value2 = 4058.8235
$ws.Range("K132").Value2 = 7602.75
$ws.Range("L132").Value2 = 12176.4705
$ws.Range("M132").Value2 = -5072.75
$ws.Range("N132").Value2 = -17236.4705

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value2 = 556.4667
$ws.Range("I100").Value2 = 544.8182
$ws.Range("J100").Value2 = 588.5
$ws.Range("K100").Value2 = 1089.6364
$ws.Range("L100").Value2 = 1177
$ws.Range("M100").Value2 = -548.6364000000001
$ws.Range("N100").Value2 = -2259
$ws.Range("H113").Value2 = 741.7222
$ws.Range("I113").Value2 = 150.66667
$ws.Range("J113").Value2 = 1923.8334
$ws.Range("K113").Value2 = 452.00001
$ws.Range("L113").Value2 = 5771.5002
$ws.Range("M113").Value2 = 1717.99999
$ws.Range("N113").Value2 = -10111.5002

Write-Output "All edits applied."
